# Add a new "note" entry to the tracking sheet:
#   B3 = 2021-03-01 (stored as serial date 44256, formatted mm-dd-yy -> numFmtId 14)
#   C3 = "Design Card Css"
# Both cells use the larger (18pt) font that matches the sheet's header font,
# and the new row keeps the same row height as the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: new note text (create the font-only style first so it gets the lower style index)
$ws.Range("C3").Value = "Design Card Css"
$ws.Range("C3").Font.Size = 18

# B3: date value with a date number format (creates the second new style)
$ws.Range("B3").Value = 44256
$ws.Range("B3").Font.Size = 18
$ws.Range("B3").NumberFormat = "mm-dd-yy"

# Keep the new row the same height as the header/data row above it
$ws.Rows("3:3").RowHeight = 32.4

# Match the saved selection state from the edit
$ws.Range("F3").Select()
